# C5-PowerPoint.pptx edit
#
# 1) Slide 6 table: switch its table style from the custom style
#    {7726B579-90EB-4E6F-A55C-407E595E90B3} to the built-in style
#    {849CAB0A-06B0-47D1-B2E8-8287A5D10389}.
# 2) Theme: the deck's theme colour scheme ("Integral") is replaced by the
#    default Office colour scheme ("Office").

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 ---------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{849CAB0A-06B0-47D1-B2E8-8287A5D10389}")
    }
}

# --- 2. Theme colours --------------------------------------------------
# Map of theme colour slot -> new (Office) RGB hex value, in the standard
# DrawingML clrScheme order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $slide6.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
